# Generate Report for Handoff
#
# The localization-status report just finished handoff generation for this
# file, so:
#   - the zh-cn / de-de status flips from "In Translation" to
#     "Ready for handoff" (Overview sheet + each language sheet's Status
#     column)
#   - the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamps advance a few dozen seconds to the new generation time
#   - the Status / datetime columns are widened so the new, longer text
#     ("Ready for handoff") fits without truncation

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn (E2) / de-de (F2) status, HO xliff date (G2) ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-07 00:54:27"

# Widen columns E and F so "Ready for handoff" is fully visible.
# (ColumnWidth snaps to the sheet's pixel grid, so this is the closest
# achievable width to the authored 17.2159881591797 target.)
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2) ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-07 00:54:22"

$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet: Status (C2) and Latest Handoff Datetime (H2) ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-07 00:54:27"

$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
